{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change (per commit \"Working on sending string using printwriter\"):\n//   - \"Bao Nguyen: ...\" contributions line gets reworded/expanded.\n//   - \"Mofan Li: ...\" contributions line gets reworded/expanded.\n//   - \"github\" -> \"Github\" in the closing remark sentence.\n\nconst body = context.document.body;\n\n// --- 1) Bao Nguyen's contributions line -------------------------------\nconst baoResults = body.search(\n  \"Bao Nguyen: Communication protocol, Possible attacks, Key modification, Message sequence, Report, Intruder simulation\",\n  { matchCase: true }\n);\nbaoResults.load(\"items\");\nawait context.sync();\n\nif (baoResults.items.length > 0) {\n  baoResults.items[0].insertText(\n    \"Bao Nguyen: Communication protocol, Key modification using poly-alphabetic cipher, Message sequence, Message encryption & decryption, Report, Intruder simulation \u2013 possible attack\",\n    \"Replace\"\n  );\n}\n\n// --- 2) Mofan Li's contributions line ----------------------------------\nconst mofanResults = body.search(\n  \"Mofan Li: Key creation, AES \u2013 Encryption, Decryption, Message integrity, Report\",\n  { matchCase: true }\n);\nmofanResults.load(\"items\");\nawait context.sync();\n\nif (mofanResults.items.length > 0) {\n  mofanResults.items[0].insertText(\n    \"Mofan Li: Key creation, Encryption & Decryption algorithm, Message integrity\",\n    \"Replace\"\n  );\n}\n\n// --- 3) \"github\" -> \"Github\" in the closing remark ----------------------\nconst githubResults = body.search(\n  \" check out our project on github at \",\n  { matchCase: true }\n);\ngithubResults.load(\"items\");\nawait context.sync();\n\nif (githubResults.items.length > 0) {\n  githubResults.items[0].insertText(\n    \" check out our project on Github at \",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n\n// --- 4) Keep the \"_GoBack\" last-edit bookmark valid ---------------------\n// The original bookmark sat inside the text that was rewritten above, so\n// re-anchor it at the end of the last paragraph we touched (mirrors how\n// Word itself relocates \"_GoBack\" to the most recent edit point).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst finalParas = body.paragraphs;\nfinalParas.load(\"items\");\nawait context.sync();\nconst mofanPara = finalParas.items[77];\nmofanPara.getRange(\"End\").insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (per commit \"Working on sending string using printwriter\"):\n#   - \"Bao Nguyen: ...\" contributions line gets reworded/expanded.\n#   - \"Mofan Li: ...\" contributions line gets reworded/expanded.\n#   - \"github\" -> \"Github\" in the closing remark sentence.\n\n$d = $word.ActiveDocument\n\n# --- 1) Bao Nguyen's contributions line --------------------------------\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\n    \"Bao Nguyen: Communication protocol, Possible attacks, Key modification, Message sequence, Report, Intruder simulation\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Bao Nguyen: Communication protocol, Key modification using poly-alphabetic cipher, Message sequence, Message encryption & decryption, Report, Intruder simulation \u2013 possible attack\",\n    1\n) | Out-Null\n\n# --- 2) Mofan Li's contributions line -----------------------------------\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n    \"Mofan Li: Key creation, AES \u2013 Encryption, Decryption, Message integrity, Report\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"Mofan Li: Key creation, Encryption & Decryption algorithm, Message integrity\",\n    1\n) | Out-Null\n\n# --- 3) \"github\" -> \"Github\" in the closing remark -----------------------\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Execute(\n    \" check out our project on github at \",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \" check out our project on Github at \",\n    1\n) | Out-Null\n"}
